$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (existing B/C shift right to C/D),
# carrying over the 75.81640625 width used by column A.
$ws.Columns.Item(2).Insert()
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(1).ColumnWidth

# New header + query for the "StatQuery" stat-bar lookup.
$ws.Cells.Item(1, 2).Value = "StatQuery"
$ws.Cells.Item(2, 2).Value = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN['']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

# Match the wrap-text style already used on A2 for the new query cell.
$ws.Cells.Item(2, 2).WrapText = $true

# Restore the view: select A2 and scroll back to the top.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("A2").Select()
